$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks (they pointed at the old TestImport* test data)
$ws.Hyperlinks.Delete()

# ---- Header row ----
$ws.Range("A1").Value = "Họ và tên"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Số điện thoại"
$ws.Range("D1").Value = "Học vị"
$ws.Range("E1").Value = "Bộ môn"
$ws.Range("F1").Value = "Khoa"
$ws.Range("G1").Value = "Chức vị"

# ---- Data rows (columns A, C-G). Column B handled separately below so the
# hyperlink + cell style + text all line up correctly. ----
# Phone numbers keep their existing text (quote-prefixed) cell format, so we
# write them with a leading apostrophe to force text and preserve the
# leading zero instead of being parsed as a number.
$ws.Range("A2").Value = "Nguyễn Văn A"
$ws.Range("C2").Value = "'0969615123"
$ws.Range("D2").Value = "Thạc sĩ"
$ws.Range("E2").Value = "Hệ thống thông tin"
$ws.Range("F2").Value = "Công nghệ thông tin"

$ws.Range("A3").Value = "Nguyễn Văn B"
$ws.Range("C3").Value = "'0969615456"
$ws.Range("D3").Value = "Thạc sĩ"
$ws.Range("F3").Value = "Công nghệ thông tin"

$ws.Range("A4").Value = "Nguyễn Văn C"
$ws.Range("C4").Value = "'0969615789"
$ws.Range("D4").Value = "Thạc sĩ"
$ws.Range("E4").Value = "Trí tuệ nhân tạo"
$ws.Range("F4").Value = "Công nghệ thông tin"

$ws.Range("A5").Value = "Nguyễn Văn D"
$ws.Range("C5").Value = "'0969615246"
$ws.Range("D5").Value = "Tiến sĩ"
$ws.Range("E5").Value = "Kỹ thuật phần mềm"
$ws.Range("F5").Value = "Công nghệ thông tin"
$ws.Range("G5").Value = "Trưởng khoa"

$ws.Range("A6").Value = "Nguyễn Văn E"
$ws.Range("C6").Value = "'0969615247"
$ws.Range("D6").Value = "Thạc sĩ"
$ws.Range("E6").Value = "Kỹ thuật phần mềm"
$ws.Range("F6").Value = "Công nghệ thông tin"

$ws.Range("A7").Value = "Nguyễn Thị G"
$ws.Range("C7").Value = "'0969615248"
$ws.Range("D7").Value = "Thạc sĩ"
$ws.Range("E7").Value = "Kỹ thuật phần mềm"
$ws.Range("F7").Value = "Công nghệ thông tin"

$ws.Range("A8").Value = "Nguyễn Thị H"
$ws.Range("C8").Value = "'0969615249"
$ws.Range("D8").Value = "Thạc sĩ"
$ws.Range("E8").Value = "Kỹ thuật phần mềm"
$ws.Range("F8").Value = "Công nghệ thông tin"

$ws.Range("A9").Value = "Nguyễn Thị I"
$ws.Range("C9").Value = "'0969615250"
$ws.Range("D9").Value = "Thạc sĩ"
$ws.Range("F9").Value = "Công nghệ thông tin"

$ws.Range("A10").Value = "Nguyễn Thị J"
$ws.Range("C10").Value = "'0969615251"
$ws.Range("D10").Value = "Thạc sĩ"
$ws.Range("F10").Value = "Công nghệ thông tin"

$ws.Range("A11").Value = "Nguyễn Thị K"
$ws.Range("C11").Value = "'0969615252"
$ws.Range("D11").Value = "Thạc sĩ"
$ws.Range("F11").Value = "Công nghệ thông tin"

# ---- Column B: emails with hyperlinks ----
# Add the hyperlink first, THEN set the cell style + text -- Hyperlinks.Add
# stamps its own (duplicate) style record on the cell, so re-applying the
# intended style/value afterwards keeps things matching the source file's
# existing style indices instead of growing new near-duplicate styles.

# Rows 2-5 use the plain Hyperlink cell style (no quote prefix)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:gv1@tlu.edu.vn")
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B2").Value = "gv1@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:gv2@tlu.edu.vn")
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B3").Value = "gv2@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:gv3@tlu.edu.vn")
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B4").Value = "gv3@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:gv4@tlu.edu.vn")
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B5").Value = "gv4@tlu.edu.vn"

# Rows 6-10 use the Hyperlink style with a text quote-prefix (as in the
# source file)
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:gv5@tlu.edu.vn")
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B6").Value = "'gv5@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:gv6@tlu.edu.vn")
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B7").Value = "'gv6@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:gv7@tlu.edu.vn")
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B8").Value = "'gv7@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:gv8@tlu.edu.vn")
$ws.Range("B9").Style = "Hyperlink"
$ws.Range("B9").Value = "'gv8@tlu.edu.vn"

$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:gv9@tlu.edu.vn")
$ws.Range("B10").Style = "Hyperlink"
$ws.Range("B10").Value = "'gv9@tlu.edu.vn"

# Row 11 keeps the same quote-prefixed Hyperlink style but (per the source
# file) has no live hyperlink relationship.
$ws.Range("B11").Style = "Hyperlink"
$ws.Range("B11").Value = "'gv10@tlu.edu.vn"

# ---- Column width for A (bestfit-like, closest achievable) ----
$ws.Columns.Item(1).ColumnWidth = 12.73

# ---- Selection ----
$ws.Range("J10").Select()
